$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 280 so the existing rows
# 280-282 shift down to become rows 282-284.
$ws.Rows.Item(280).Resize(2).Insert()

# New (current week) data for row 280 - "Primera"
$ws.Range("A280").Value = 3
$ws.Range("B280").Value = "Femacal de La Calera"
$ws.Range("C280").Value = "Coquimbo"
$ws.Range("D280").Value = 45121
$ws.Range("D280").Style = $ws.Range("D282").Style
$ws.Range("E280").Value = 5
$ws.Range("F280").Value = "Fruta"
$ws.Range("G280").Value = 100107
$ws.Range("H280").Value = "Otros"
$ws.Range("I280").Value = 100107002
$ws.Range("J280").Value = "Chirimoya"
$ws.Range("K280").Value = "Cultivar IV Región"
$ws.Range("L280").Value = "Primera"
$ws.Range("M280").Value = 36
$ws.Range("N280").Value = 30000
$ws.Range("O280").Value = 30000
$ws.Range("P280").Value = 30000
$ws.Range("Q280").Value = "$/bandeja 10 kilos"
$ws.Range("R280").Value = "Provincia del Elquí"
$ws.Range("S280").Value = 3000
$ws.Range("T280").Value = 10

# New (current week) data for row 281 - "Segunda"
$ws.Range("A281").Value = 3
$ws.Range("B281").Value = "Femacal de La Calera"
$ws.Range("C281").Value = "Coquimbo"
$ws.Range("D281").Value = 45121
$ws.Range("D281").Style = $ws.Range("D282").Style
$ws.Range("E281").Value = 5
$ws.Range("F281").Value = "Fruta"
$ws.Range("G281").Value = 100107
$ws.Range("H281").Value = "Otros"
$ws.Range("I281").Value = 100107002
$ws.Range("J281").Value = "Chirimoya"
$ws.Range("K281").Value = "Cultivar IV Región"
$ws.Range("L281").Value = "Segunda"
$ws.Range("M281").Value = 30
$ws.Range("N281").Value = 27000
$ws.Range("O281").Value = 27000
$ws.Range("P281").Value = 27000
$ws.Range("Q281").Value = "$/bandeja 10 kilos"
$ws.Range("R281").Value = "Provincia del Elquí"
$ws.Range("S281").Value = 2700
$ws.Range("T281").Value = 10
